$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 (Order = 7, Primary/Foreign Keys...) is missing the "Details" note
# that appears on the very similar row 9. Fill it in to match.
$ws.Range("F8").Value = "Slides no demo"
